$d = $word.ActiveDocument

# ------------------------------------------------------------------
# For each "github" occurrence that needs to become "GitHub" we
# toggle a character-level property (Bold) on, overwrite the two-
# word-boundary text, then toggle it back off. That forces Word to
# split the run at the match boundaries (mirroring what real Word
# does when AutoCorrect / Find-Replace touches only part of a run),
# giving us separate <w:r> runs for "GitHub" instead of merging it
# back into the surrounding text.
# ------------------------------------------------------------------

# ---- Part 1: "Make sure the author is: <github username> <github user email>" ----
$para1 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Make sure the author is:*") {
        $para1 = $p
    }
}
$pstart = $para1.Range.Start

# First occurrence of "github"
$full = $para1.Range.Text
$idx1 = $full.IndexOf("github")
$r1 = $d.Range($pstart + $idx1, $pstart + $idx1 + 6)
$r1.Bold = 1
$r1.Text = "GitHub"
$r1b = $d.Range($pstart + $idx1, $pstart + $idx1 + 6)
$r1b.Bold = 0

# Second occurrence of "github" (recompute text/offset after first edit)
$full2 = $para1.Range.Text
$idx2 = $full2.IndexOf("github", $idx1 + 6)
$r2 = $d.Range($pstart + $idx2, $pstart + $idx2 + 6)
$r2.Bold = 1
$r2.Text = "GitHub"
$r2b = $d.Range($pstart + $idx2, $pstart + $idx2 + 6)
$r2b.Bold = 0

# ---- Part 2: "Enter the github username and password and press OK." ----
$para2 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Enter the github username*") {
        $para2 = $p
    }
}
$pstart2 = $para2.Range.Start
$fullE = $para2.Range.Text
$idxE = $fullE.IndexOf("github")
$rE = $d.Range($pstart2 + $idxE, $pstart2 + $idxE + 6)
$rE.Bold = 1
$rE.Text = "GitHub"
$rEb = $d.Range($pstart2 + $idxE, $pstart2 + $idxE + 6)
$rEb.Bold = 0

# ---- Part 3: move the "_GoBack" bookmark from the "Team > Pull" ----
# paragraph to right after the newly-capitalized "GitHub" in Part 2.
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

$bmPos = $pstart2 + $idxE + 6
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
